$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row above row 48, shifting existing rows (48:203) down to (49:204)
$ws.Rows("48:48").Insert()

# Populate the newly inserted row with the new transaction entry
$ws.Range("R48").Value = "fake messages suspicious"
$ws.Range("S48").Value = "2024-09-24 12:23:30"
